$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet right after "总计" (i.e. before "2022-Q1").
#    Clone the "2022-Q1" sheet (same column layout/headers/styles) into that
#    slot, rename it, then overwrite its data rows with the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Copy($null, $wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Row 2 - 970042 (overwrite the cloned 2022-Q1 row)
$newSheet.Range("B2").Value = "'970042"
$newSheet.Range("C2").Value = "国海量化优选一年持有股票C"
$newSheet.Range("D2").Value = "'7.16"
$newSheet.Range("E2").Value = "'87.31"
$newSheet.Range("F2").Value = "'0.36"
$newSheet.Range("G2").Value = "'0.0258"
$newSheet.Range("H2").Value = 1

# Row 3 - 970041
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'970041"
$newSheet.Range("C3").Value = "国海量化优选一年持有股票A"
$newSheet.Range("D3").Value = "'0.63"
$newSheet.Range("E3").Value = "'87.31"
$newSheet.Range("F3").Value = "'0.36"
$newSheet.Range("G3").Value = "'0.0023"
$newSheet.Range("H3").Value = 1

# Drop the "quote prefix" style that typing a leading apostrophe adds, while
# keeping the text type/value intact - matches the plain (unstyled) text
# cells used on every other quarter sheet.
$newSheet.Range("B2:G3").ClearFormats()

# Re-apply the header / index-column look-and-feel (bold, bordered, centered)
# used throughout the workbook by copying the format from an existing sheet.
$wb.Worksheets.Item("2022-Q1").Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$wb.Worksheets.Item("2022-Q1").Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("总计").Rows.Item(2).Insert()

# The inserted row inherits the header row's bold/bordered look via
# "format like row above" - strip that back to plain formatting.
$wb.Worksheets.Item("总计").Range("B2:D2").ClearFormats()

# Give the new A2 index cell the same bold/bordered/centered style as the
# rest of column A.
$wb.Worksheets.Item("总计").Range("A3").Copy()
$wb.Worksheets.Item("总计").Range("A2").PasteSpecial(-4122)

# New row 2 values (2022-Q3 summary).
$wb.Worksheets.Item("总计").Range("A2").Value = 0
$wb.Worksheets.Item("总计").Range("B2").Value = "2022-Q3"
$wb.Worksheets.Item("总计").Range("C2").Value = 2
$wb.Worksheets.Item("总计").Range("D2").Value = 0.03

# Column A is a running index (0,1,2,...) independent of which quarter sits
# in the row, so renumber it after the insert shifted the old rows down.
$wb.Worksheets.Item("总计").Range("A3").Value = 1
$wb.Worksheets.Item("总计").Range("A4").Value = 2
$wb.Worksheets.Item("总计").Range("A5").Value = 3
$wb.Worksheets.Item("总计").Range("A6").Value = 4
$wb.Worksheets.Item("总计").Range("A7").Value = 5
$wb.Worksheets.Item("总计").Range("A8").Value = 6
